# Generate Report for Handback
#
# Marks the two handed-off files (281112d5-...md and 72dfbeb8-...md) as
# handed back, for both the zh-cn and de-de localization targets:
#   - Overview / per-locale Status columns switch from "Ready for handoff"
#     to "Handed back: in sync with en-US"
#   - Latest Target File (I) gets the handoff source file as a hyperlink
#   - Latest Handback File (J) gets the generated xlf file name
#   - Latest Handback DateTime (K) gets the handback timestamp
#   - Relevant columns are widened so the new content is readable

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$md1 = "281112d5-af5a-4b77-ac95-317023dcf2e9.md"
$md2 = "72dfbeb8-593a-400d-a8d4-77237649c29b.md"
$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2a78fabeb79efd4b344c8702d0df0ac6f1118e4/e2e/281112d5-af5a-4b77-ac95-317023dcf2e9.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2a78fabeb79efd4b344c8702d0df0ac6f1118e4/e2e/72dfbeb8-593a-400d-a8d4-77237649c29b.md"

function Update-LocaleSheet {
    param($SheetName, $Xlf1, $Xlf2, $HandbackTime)

    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    $ws.Range("J2").Value = $Xlf1
    $ws.Range("J3").Value = $Xlf2

    $ws.Range("K2").Value = $HandbackTime
    $ws.Range("K3").Value = $HandbackTime

    $ws.Hyperlinks.Add($ws.Range("I2"), $url1, "", "", $md1)
    $ws.Hyperlinks.Add($ws.Range("I3"), $url2, "", "", $md2)

    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

Update-LocaleSheet "zh-cn" "281112d5-af5a-4b77-ac95-317023dcf2e9.87c83b646adcf3a265a49c0c022b9aa4c78b8642.zh-cn.xlf" "72dfbeb8-593a-400d-a8d4-77237649c29b.c6f187c302caabd680fd62d118c073449a94aea0.zh-cn.xlf" "2016-08-15 22:44:42"

Update-LocaleSheet "de-de" "281112d5-af5a-4b77-ac95-317023dcf2e9.87c83b646adcf3a265a49c0c022b9aa4c78b8642.de-de.xlf" "72dfbeb8-593a-400d-a8d4-77237649c29b.c6f187c302caabd680fd62d118c073449a94aea0.de-de.xlf" "2016-08-15 22:44:49"
